$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.336.82'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '1.710.77'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5292'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06681'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2660'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.82'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07688'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.504'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.41%  '
$ws.Range("D13").Value = '1.946.12'
$ws.Range("E13").Value = '  -0.84%  '
$ws.Range("D14").Value = '1.707.99'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5853'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '0.0₅8212'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D18").Value = '27.357.25'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.642'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.014'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.686'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1209'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.236'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05336'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.292'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.468'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.426'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.71%  '
$ws.Range("E34").Value = '  -1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.874'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9527'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.394'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5851'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").Value = '1.145.85'
$ws.Range("E39").Value = '  +8.40%  '
$ws.Range("E40").Value = '  -0.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.788'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.96%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8387'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").Value = '1.853.21'
$ws.Range("E45").Value = '  -0.84%  '
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.52'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4566'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.36%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.098'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05206'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.69%  '
